# Automated_SmokeTest_Result.xlsx update
# New smoke-test run: refreshed timestamps, refreshed 'TABLE'/'FORM'/'REPORT' page
# verification values, and a couple of description/title tweaks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("C1").Value = "20/05/2024 08:23:AM"

# --- File submission section ------------------------------------------
$ws.Range("B22").Value = "9_INVALID_ROWS_FORM_S.txt"
$ws.Range("C22").Value = "9_INVALID_ROWS_FORM_S.txt"

$ws.Range("B24").Value = "Invalid Rows"
$ws.Range("C24").Value = "Invalid Rows"

$ws.Range("B26").Value = "05/20/2024 08:24 AM"
$ws.Range("C26").Value = "05/20/2024 08:24 AM"

# --- 'TABLE' page section (was 'DASHBOARD') ----------------------------
$ws.Range("A28").Value = "Verifying 'TABLE' Page's First Record & MicroStrategy Title"

$ws.Range("A30").Value = "Table ID"
$ws.Range("B30").Value = "UTBL01"
$ws.Range("C30").Value = "UTBL01"

$ws.Range("A32").Value = "Table Type"
$ws.Range("B32").Value = "PART B"
$ws.Range("C32").Value = "PART B"

$ws.Range("A34").Value = "Table Name"
$ws.Range("B34").Value = "WORKLOAD OPERATIONS: CLAIMS RECEIVED, PROCESSED, AND PENDING DATA"
$ws.Range("C34").Value = "WORKLOAD OPERATIONS: CLAIMS RECEIVED, PROCESSED, AND PENDING DATA"

$ws.Range("B36").Value = "(UTBL01 - PART B) WORKLOAD OPERATIONS: CLAIMS RECEIVED, PROCESSED, AND PENDING DATA. MicroStrategy"
$ws.Range("C36").Value = "(UTBL01 - PART B) WORKLOAD OPERATIONS: CLAIMS RECEIVED, PROCESSED, AND PENDING DATA. MicroStrategy"
$ws.Range("D36").Value = "Fail"

# --- 'FORM' page section -------------------------------------------------
$ws.Range("B41").Value = "FORM 7"
$ws.Range("C41").Value = "FORM 7"
# D41 keeps "Pass" but picks up a fresh (automatic-color) font, as in the source edit
$ws.Range("D41").Font.Color = 0

$ws.Range("C45").Value = "APPEALS ACTIVITY (CMS-2592)"
$ws.Range("D45").ClearContents()

# --- 'REPORT' page section -----------------------------------------------
$ws.Range("A50").Value = "Verifying 'REPORT' Page's First Record "

$ws.Range("B52").Value = "URPT01"
$ws.Range("C52").Value = "URPT01"

$ws.Range("B54").Value = "CONTRACTOR MAPPING"
$ws.Range("C54").Value = "CONTRACTOR MAPPING"
# D54 keeps "Pass" but picks up the same fresh font as D41
$ws.Range("D54").Font.Color = 0

$ws.Rows.Item(56).RowHeight = 18.6
$ws.Range("B56").Value = "This report lists all CROWD Contractor Details including their roles."
$ws.Range("C56").Value = "This report lists all CROWD Contractor Details including their roles."

# --- 'Resources' page section ---------------------------------------------
$ws.Range("A58").Value = "Verifying 'Resources' Page's First Record"

# --- 'NEWS' page section ----------------------------------------------------
$ws.Range("A67").Value = "Verifying 'NEWS' Page For Year 2020"

$ws.Range("B69").Value = "Attn: ALL CROWD Users | October 1, 2024"
$ws.Range("C69").Value = "Attn: ALL CROWD Users | October 1, 2024"

$ws.Rows.Item(70).RowHeight = 33.6
$ws.Range("B70").Value = "This is a placeholder. This will be updated closer to the Go Live Date."
$ws.Range("B70").Style = "Normal"
$ws.Range("C70").Value = "This is a placeholder. This will be updated closer to the Go Live Date."

# --- Column widths (closest achievable snap on this engine's width grid) ---
$ws.Columns.Item(2).ColumnWidth = 95.65
$ws.Columns.Item(3).ColumnWidth = 100.8

# --- Selection / view ------------------------------------------------------
$ws.Range("F55").Select()
